# Apply the "Conclusion" slide edit:
#  - merge the "ANN Bal. model" paragraph into the preceding paragraph's
#    sentence (appending a period) and drop the now-redundant paragraph
#  - shrink the content placeholder to match the shorter text
#  - drop the orphaned click-animation that targeted the removed paragraph

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)

# Rewrite the body text: keep paragraph 1 untouched, replace paragraph 2's
# text with the period-terminated sentence, and remove paragraph 3 entirely.
$tr = $shp.TextFrame.TextRange
$tr.Text = "Balanced, cost-sensitive CART (CART Cost Bal.) model" + [char]13 + "Balanced artificial neural network (ANN Bal.) model."

# Shrink the placeholder now that it only holds two lines of text.
$shp.Height = 81.9698725

# The third paragraph's click-triggered entrance animation is now orphaned
# (it targeted paragraph index 2, which no longer exists) -- remove it.
$seq = $s.TimeLine.MainSequence
$seq.Item($seq.Count).Delete()
